# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E77) previously listed the arrears periods
# in descending order (2301 .. 1712). The data refresh re-lists them in
# ascending order (1712 .. 2301), and the single irregular "Valor Mora"
# amount (24591, everything else is 29509) now belongs to the last period
# (2301 / row 77) instead of the first (1712 / row 16).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301"
)

$firstRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $ws.Cells.Item($firstRow + $i, 5).Value = $periods[$i]
}

# The odd "Valor Mora" value now sits on the last period row instead of the first.
$ws.Cells.Item(16, 6).Value = 29509
$ws.Cells.Item(77, 6).Value = 24591
